# Revert the rule for rev 0 so "0" revisions resolve to "0" instead of blank,
# for compatibility with the MWS software.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("components request")

# Update the formulas in C2 and C3 to use 1 (not blank) when the last char of B is "0"
$ws.Range("C2").Formula = '=RIGHT(B2,IF(RIGHT(B2,1)="0",1,2))'
$ws.Range("C3").Formula = '=RIGHT(B3,IF(RIGHT(B3,1)="0",1,2))'

# Update the active selection on this sheet to C3
$ws.Activate()
$ws.Range("C3").Select()
